# Applies "elapsed time y cpu" edit: adds two new columns (G: Elapsed Time, H: CPU)
# with header style matching the existing header row, and fills in updated
# B/C/D metric values plus new G/H values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (copy style from an existing header cell, e.g. F1) ---
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("B2").Value = 0.4964150075448933
$ws.Range("C2").Value = 0.9854189908374392
$ws.Range("D2").Value = 0.5197943089546899
$ws.Range("G2").Value = 0.4794827245333484
$ws.Range("H2").Value = 0.996

# --- Row 3 ---
$ws.Range("B3").Value = 0.09433148089300998
$ws.Range("C3").Value = 0.9986963133924546
$ws.Range("D3").Value = 0.2490681084032768
$ws.Range("G3").Value = 0.4794827245333484
$ws.Range("H3").Value = 0.996

# --- Row 4 ---
$ws.Range("B4").Value = 0.03356298177852739
$ws.Range("C4").Value = 0.9996540735033265
$ws.Range("D4").Value = 0.132108406813797
$ws.Range("G4").Value = 0.4794827245333484
$ws.Range("H4").Value = 0.996

# --- Row 5 ---
$ws.Range("B5").Value = 0.09838959378682038
$ws.Range("C5").Value = 0.9994139212625992
$ws.Range("D5").Value = 0.2210435799373883
$ws.Range("G5").Value = 0.4794827245333484
$ws.Range("H5").Value = 0.996
